# Updated symbol list on Tue Dec 20 18:35:13 UTC 2022 with GitHub Actions
#
# Applies the per-cell updates to the "cryptos" sheet: refreshed prices in
# column D, a few Volume(1h) label tweaks in column E, and a handful of rows
# (14-18 and 25-26) whose Coin/Link/Price/Volume values shifted because the
# underlying coin ranking list was re-ordered.
#
# Price values in column D are stored as text (not numbers) in the workbook,
# so we prefix them with a leading apostrophe when assigning via .Value to
# force Excel to keep them as text instead of auto-converting to numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D price refresh (rows 2-12) ---
$ws.Range("D2").Value  = "'249.13"
$ws.Range("D3").Value  = "'22.85"
$ws.Range("D4").Value  = "'5.426"
$ws.Range("D5").Value  = "'0.05631"
$ws.Range("D6").Value  = "'3.420"
$ws.Range("D7").Value  = "'6.356"
$ws.Range("D8").Value  = "'0.8131"
$ws.Range("D9").Value  = "'0.9177"
$ws.Range("D10").Value = "'0.1437"
$ws.Range("D11").Value = "'0.07529"
$ws.Range("D12").Value = "'0.03103"

# --- Rows 14-18: coin list re-ordered (ProBitToken inserted before BitMartToken) ---
$ws.Range("B14").Value = "ProBitToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D14").Value = "'0.1297"
$ws.Range("E14").Value = "13ProBitTokenPROB"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09324"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.568"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001581"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04768"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- Column D price refresh (rows 19-24) ---
$ws.Range("D19").Value = "'0.006385"
$ws.Range("D20").Value = "'0.004997"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.699"
$ws.Range("D24").Value = "'2.178"

# --- Rows 25-26: coin list re-ordered (One moved up, BitpandaEcosystemToken follows) ---
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D25").Value = "'0.01152"
$ws.Range("E25").Value = "24OneONEBestin24h"

$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3300"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"

# --- Misc label/price tweaks (rows 27-50) ---
$ws.Range("E27").Value = "26AAXTokenAABWorstin24h"
$ws.Range("D28").Value = "'0.0003032"

$ws.Range("D40").Value = "'0.04019"
$ws.Range("D41").Value = "'0.006801"
$ws.Range("D42").Value = "'0.1068"
$ws.Range("D43").Value = "'0.002713"
$ws.Range("D44").Value = "'0.007500"
$ws.Range("D45").Value = "'0.00005467"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.5000"
$ws.Range("D48").Value = "'0.2407"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.01010"

Write-Host "Applied cryptos sheet updates"
